$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily food-log entries for 22-04-2018 through 29-04-2018, continuing
# on from the last existing row (111 = 21-04-2018).
$data = @(
    @("22-04-2018", 2058, 61, 21, 230, 3285, 147, 3000),
    @("23-04-2018", 1891, 52, 30, 174, 4042, 181, 1000),
    @("24-04-2018", 1658, 46, 16, 145, 2632, 155, 3000),
    @("25-04-2018", 1571, 45, 11, 126, 3182, 149, 2500),
    @("26-04-2018", 2050, 111, 9, 86, 2649, 166, 3000),
    @("27-04-2018", 2010, 77, 13, 191, 4375, 128, 1750),
    @("28-04-2018", 1832, 62, 10, 158, 2212, 133, 1500),
    @("29-04-2018", 2236, 104, 23, 183, 4231, 148, 3250)
)

$lastRow = 111
$startRow = $lastRow + 1
$endRow = $startRow + $data.Length - 1

# Seed the new rows by copying the formatting (and formulas) of the last
# existing row down, the same way a user would fill the log forward.
$srcRow = $ws.Range("A$lastRow" + ":J$lastRow")
for ($row = $startRow; $row -le $endRow; $row++) {
    $destRow = $ws.Range("A$row" + ":J$row")
    $srcRow.Copy($destRow)
}

# Now fill in the real values and formulas for each new day.
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]

    $ws.Cells.Item($row, 9).Formula = "=IF(H$row>=2200,""Yes"",""No"")"
    $ws.Cells.Item($row, 10).Formula = "=IF(B$row<=1800,""Yes"",""No"")"
}

# Match the author's final view/selection state.
$ws.Range("L113").Select()
